# Update "Income Statements" workbook: drop the oldest fiscal-year column
# (1396/12) and its publish date, shift the remaining four periods one
# column to the left, and append a new rightmost period (1401/12) with a
# refreshed publish date and recomputed financial figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: period headers -------------------------------------------------
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish dates ----------------------------------------------------
$ws.Range("D9").Value = "1399-04-21 (13)"
$ws.Range("E9").Value = "1400-02-30 (7)"
$ws.Range("F9").Value = "1401-04-18 (7)"
$ws.Range("G9").Value = "1402-02-23 (8)"

# H9 ("1402-02-23") looks like an ISO date, so Excel's smart-typing would
# otherwise silently convert it to a date serial. Force text storage, then
# restore H9's visual formatting (border/fill/font) from its neighbour so
# it keeps matching the rest of the row instead of picking up a bespoke
# "@" number-format style.
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "1402-02-23"
$ws.Range("G9").Copy()
$ws.Range("H9").PasteSpecial(-4122)

# --- Row 11: فروش (Sales) ---------------------------------------------------
$ws.Range("D11").Value = 61358
$ws.Range("E11").Value = 78969
$ws.Range("F11").Value = 79805
$ws.Range("G11").Value = 87595
$ws.Range("H11").Value = 106811

# --- Row 12: بهای تمام شده کالای فروش رفته (Cost of goods sold) -------------
$ws.Range("D12").Value = -39032
$ws.Range("E12").Value = -50762
$ws.Range("F12").Value = -56156
$ws.Range("G12").Value = -67298
$ws.Range("H12").Value = -74576

# --- Row 13: سود (زیان) ناخالص (Gross profit) --------------------------------
$ws.Range("D13").Value = 22326
$ws.Range("E13").Value = 28206
$ws.Range("F13").Value = 23649
$ws.Range("G13").Value = 20297
$ws.Range("H13").Value = 32234

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses) ---------------
$ws.Range("D14").Value = -5967
$ws.Range("E14").Value = -5627
$ws.Range("F14").Value = -5966
$ws.Range("G14").Value = -8568
$ws.Range("H14").Value = -9452

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی --------------------------
$ws.Range("D16").Value = 115
$ws.Range("E16").Value = 18
$ws.Range("F16").Value = 841
$ws.Range("G16").Value = 235
$ws.Range("H16").Value = -120

# --- Row 17: سود (زیان) عملیاتی (Operating profit) ---------------------------
$ws.Range("D17").Value = 16474
$ws.Range("E17").Value = 22597
$ws.Range("F17").Value = 18524
$ws.Range("G17").Value = 11964
$ws.Range("H17").Value = 22663

# --- Row 18: هزینه های مالی (Financial expenses) -----------------------------
$ws.Range("D18").Value = -3432
$ws.Range("E18").Value = -4961
$ws.Range("F18").Value = -4222
$ws.Range("G18").Value = -5649
$ws.Range("H18").Value = -7061

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی -------------------------
$ws.Range("D19").Value = 1975
$ws.Range("E19").Value = 3212
$ws.Range("F19").Value = 4824
$ws.Range("G19").Value = 4480
$ws.Range("H19").Value = 3266

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات ---------------
$ws.Range("D20").Value = 15018
$ws.Range("E20").Value = 20848
$ws.Range("F20").Value = 19126
$ws.Range("G20").Value = 10794
$ws.Range("H20").Value = 18868

# --- Row 21: مالیات (Tax) -----------------------------------------------------
$ws.Range("D21").Value = -1354
$ws.Range("E21").Value = -1127
$ws.Range("F21").Value = -688
$ws.Range("G21").Value = -486
$ws.Range("H21").Value = -147

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم ------------------------------
$ws.Range("D22").Value = 13663
$ws.Range("E22").Value = 19721
$ws.Range("F22").Value = 18438
$ws.Range("G22").Value = 10309
$ws.Range("H22").Value = 18721

# --- Row 24: سود (زیان) خالص (Net profit) -------------------------------------
$ws.Range("D24").Value = 13663
$ws.Range("E24").Value = 19721
$ws.Range("F24").Value = 18438
$ws.Range("G24").Value = 10309
$ws.Range("H24").Value = 18721

# --- Row 26: سرمایه (Capital) --------------------------------------------------
$ws.Range("D26").Value = 5931
$ws.Range("E26").Value = 7795
$ws.Range("F26").Value = 4423
$ws.Range("G26").Value = 7580
$ws.Range("H26").Value = 11335
